$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11 - this pushes the existing "cell4"/"angle4"/
# "cell5"/"angle5" blocks (old rows 11,12,14,15) down to rows 12,13,15,16,
# leaving a fresh blank row 10 (row 11 stays blank too, preserving the
# existing cell4/blank/angle4/blank grouping pattern).
$ws.Rows.Item(11).Insert()

# Populate the new row 10 with the "size3" label and its raw measurements.
$ws.Cells.Item(10, 1).Value = "size3"
$ws.Cells.Item(10, 2).Value = 108
$ws.Cells.Item(10, 3).Value = 107
$ws.Cells.Item(10, 4).Value = 103
$ws.Cells.Item(10, 5).Value = 102
$ws.Cells.Item(10, 6).Value = 103
$ws.Cells.Item(10, 7).Value = 103
$ws.Cells.Item(10, 8).Value = 109
$ws.Cells.Item(10, 9).Value = 131
$ws.Cells.Item(10, 10).Value = 263
$ws.Cells.Item(10, 11).Value = 230

# Move the active selection to K10, matching the saved view state.
[void]$ws.Range("K10").Select()
